# delete samples for running faster
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 9 and 10 first (T1043, T112), then rows 2 and 3 (T1013, T1015),
# working bottom-up so row numbers of not-yet-deleted rows stay valid.
$ws.Range("A9:C10").EntireRow.Delete()
$ws.Range("A2:C3").EntireRow.Delete()
